$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column headers (row 2): M2 = "ТК_оригинал", O2 = "на момент выгрузки в элжуре"
$ws.Range("M2").Value = "ТК_оригинал"
$ws.Range("O2").Value = "на момент выгрузки в элжуре"

# Values for column O (rows 4-30) mirror column M for each student
$mValues = @(5,5,4,4,5,5,3,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,3,5)
for ($i = 0; $i -lt $mValues.Count; $i++) {
    $row = 4 + $i
    $ws.Range("O$row").Value = $mValues[$i]
}

# Apply the bordered / centered / wrapped style to O4:O30 by copying the
# existing header-like border format (thick box border, style index 1)
# and stripping the bold font so it matches a plain-body cell.
$ws.Range("C2").Copy()
$ws.Range("O4:O30").PasteSpecial(-4122)
$ws.Range("O4:O30").Font.Bold = $false

# re-apply values (PasteSpecial of formats only shouldn't touch values, but
# make sure values are correct after formatting)
for ($i = 0; $i -lt $mValues.Count; $i++) {
    $row = 4 + $i
    $ws.Range("O$row").Value = $mValues[$i]
}

# Column P: difference formula O-M, rows 4 (plain) then 5:30 (shared formula)
$ws.Range("P4").Formula = "=O4-M4"
$ws.Range("P5:P30").Formula = "=O5-M5"

# Update the worksheet view: rescroll, rezoom, reselect active cell
$ws.Activate()
$ws.Range("C4").Select()
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("C4").Select()
$excel.ActiveWindow.FreezePanes = $true
$excel.ActiveWindow.Zoom = 100
$ws.Range("M28").Select()
